$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# The workbook was refreshed a day later (03-Nov-2025 -> 04-Nov-2025):
# each training row's "PERIOD TO EXPIRE" (col H) ticks down by one day
# and "LAST UPDATE" (col I) moves to the new date, for every data row
# (3 through 29).
for ($row = 3; $row -le 29; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $currentPeriod = $hCell.Value2
    $hCell.Value = $currentPeriod - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    # Leading apostrophe forces this to stay literal text (like the
    # original "03-Nov-2025" inline string) instead of being
    # auto-parsed into a date serial value.
    $iCell.Value = "'04-Nov-2025"
}
